$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$range = $ws.Range("A2:A27")
$range.NumberFormat = "@"
$range.Value = "08.15.19"

$ws.Range("A3:A27").Select()
